$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.825.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.429.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.76%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.439.88"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.546"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.123"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.439"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.044.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.48%  "
$ws.Range("E14").Value = "  -2.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000193"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.973.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.421.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "389.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.11%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.97%  "
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.540"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000119"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +24.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.177"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.47"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0779"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.98%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.916.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.45%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0318"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.63%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.772"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.12%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.53%  "
$ws.Range("E48").Value = "  +3.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +22.48%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.43%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.849"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.90%  "
